$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Session"

# Delete the last data row (row 42) entirely, shifting cells up
$ws.Rows.Item(42).Delete()
